# The workbook records UserName / Password / Login Status test data on
# Sheet1. The B4 password cell ("leo_12345") was left behind as stale/
# duplicate test data and is being cleared out (handling the "null cell
# value" case referenced in the commit message) while the rest of the
# grid is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the stale password value in B4 - it becomes a genuinely empty
# (null) cell rather than holding leftover text.
$ws.Range("B4").Value = $null

# Auto-sized column widths for the UserName / Password columns so the
# (now longer) e-mail addresses and the login-status text are fully
# visible, matching the widths Excel would compute for this content.
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(2).ColumnWidth = 10.666666666666666

# Leave the selection where the author last left it when saving.
$ws.Range("C10").Select() | Out-Null
